# Apply updated crypto price/volume data (and one inserted "LEO" row)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.971.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.882.64'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.87'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4992'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3847'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09186'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -5.65%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.61'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.331'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.71'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.878.47'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.289'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001104'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.78%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.26'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06627'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.98'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.166'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.032.34'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.37'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.307'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.376'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.097.56'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.542'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.05%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '158.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.28%  '
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.75'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.47%  '
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '126.47'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1055'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.56%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.065'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.72%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.582'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.35%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.584'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.54%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.380'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.28%  '
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06568'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.35%  '
$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02397'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2183'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.284'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.85%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.203'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.62%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.58'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.12%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6410'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.24%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'InternetComputer(DFINITY)'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.931'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.08%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.33'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.11%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6015'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.57%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.292'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.72%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.669'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.55%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.987'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.201'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.74%  '
